$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - headers (columns shift: a scooter/motorbike/2-axis split is
# inserted, bicycle moves right, bus ends up in the new last column)
$ws.Range("B2").Value = "sitp"
$ws.Range("C2").Value = "scooter"
$ws.Range("D2").Value = "motorbike"
$ws.Range("E2").Value = "car"
$ws.Range("F2").Value = "bicycle"
$ws.Range("G2").Value = "2-axis"
$ws.Range("H2").Value = "bus"

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 294
$ws.Range("E3").Value = 531
$ws.Range("F3").Value = 28
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 47

# Row 4
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 346
$ws.Range("E4").Value = 470
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 30

# Row 5
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 651
$ws.Range("E5").Value = 1040
$ws.Range("F5").Value = 94
$ws.Range("G5").Value = 16
$ws.Range("H5").Value = 47

# Row 6 - the movement id also changes (5 -> 6); keep it textual like the
# other id cells in column A (leading apostrophe forces text, matching
# how the rest of the column is already stored)
$ws.Range("A6").Value = "'6"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 152
$ws.Range("E6").Value = 352
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 22
$ws.Range("H6").Value = 21

# Row 7 - Grand Total: extend the existing SUM formulas into the two new
# columns created by the header split
$ws.Range("G7").Formula = "=SUM(G3:G6)"
$ws.Range("H7").Formula = "=SUM(H3:H6)"
